$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

$newValues = @(
    131.300928,
    119.22372,
    117.056016,
    117.210852,
    117.36568800000001,
    119.843064,
    117.210852,
    118.294704,
    117.210852,
    117.210852,
    118.294704,
    117.36568800000001,
    117.056016,
    116.74634399999999,
    117.52052399999999,
    118.604376,
    116.901179999999,
    113.030279999999,
    117.675359999999,
    117.056016,
    117.210852,
    117.985032,
    117.830196,
    119.843064,
    117.985032,
    117.675359999999,
    117.675359999999,
    117.36568800000001,
    116.901179999999,
    117.210852
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}

# Update the view state: scroll position and active selection
$ws.Range("D20").Select()
$activeWindow = $excel.ActiveWindow
$activeWindow.ScrollRow = 13
$activeWindow.ScrollColumn = 1
